# Restore edit: update the "R30" rule's lower bound (cell C10) from 18 to 1.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("C10").Value = 1
